$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2671
$ws.Range("F3").Value = 1023
$ws.Range("F4").Value = 19192
$ws.Range("F6").Value = 2144
$ws.Range("F10").Value = 665
$ws.Range("F12").Value = 230
$ws.Range("F14").Value = 344
$ws.Range("F16").Value = 242
$ws.Range("F17").Value = 166
$ws.Range("F18").Value = 159
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 118
$ws.Range("F11").Value = 9
$ws.Range("F13").Value = 81
$ws.Range("F15").Value = 52
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5957
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 5957
$ws.Range("F7").Value = 2671
$ws.Range("F8").Value = 1023
$ws.Range("F9").Value = 19192
$ws.Range("F16").Value = 2144
$ws.Range("F18").Value = 118
$ws.Range("F21").Value = 665
$ws.Range("F23").Value = 230
$ws.Range("F27").Value = 9
$ws.Range("F28").Value = 344
$ws.Range("F31").Value = 242
$ws.Range("F32").Value = 81
$ws.Range("F33").Value = 166
$ws.Range("F35").Value = 159
$ws.Range("F36").Value = 52
